$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 134.647784
$ws.Range("H2").Value = 403.943352
$ws.Range("I2").Value = 0.2617460968718581
$ws.Range("J2").Value = 0.2617460968718581
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 47.57896333333333
$ws.Range("N2").Value = 142.73689
$ws.Range("O2").Value = 0.450188452948237
$ws.Range("P2").Value = 0.4501884529482371
$ws.Range("Q2").Value = 6406.401977850586
$ws.Range("R2").Value = 57657.61780065527
$ws.Range("S2").Value = 0.1178350704159812
$ws.Range("T2").Value = 0.1178350704159812
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 134.647784
$ws.Range("H3").Value = 403.943352
$ws.Range("I3").Value = 0.2617460968718581
$ws.Range("J3").Value = 0.2617460968718581
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.829723666666666
$ws.Range("N3").Value = 29.489171
$ws.Range("O3").Value = 0.09300808131111737
$ws.Range("P3").Value = 0.09300808131111739
$ws.Range("Q3").Value = 1323.550509049021
$ws.Range("R3").Value = 11911.95458144119
$ws.Range("S3").Value = 0.02434450226072539
$ws.Range("T3").Value = 0.02434450226072539
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 134.647784
$ws.Range("H4").Value = 403.943352
$ws.Range("I4").Value = 0.2617460968718581
$ws.Range("J4").Value = 0.2617460968718581
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 10.933664
$ws.Range("N4").Value = 32.800992
$ws.Range("O4").Value = 0.1034534789405002
$ws.Range("P4").Value = 0.1034534789405003
$ws.Range("Q4").Value = 1472.193628600576
$ws.Range("R4").Value = 13249.74265740518
$ws.Range("S4").Value = 0.02707854432049091
$ws.Range("T4").Value = 0.02707854432049091
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 134.647784
$ws.Range("H5").Value = 403.943352
$ws.Range("I5").Value = 0.2617460968718581
$ws.Range("J5").Value = 0.2617460968718581
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 37.34441866666666
$ws.Range("N5").Value = 112.033256
$ws.Range("O5").Value = 0.3533499868001453
$ws.Range("P5").Value = 0.3533499868001453
$ws.Range("Q5").Value = 5028.343218234901
$ws.Range("R5").Value = 45255.08896411411
$ws.Range("S5").Value = 0.09248797987466062
$ws.Range("T5").Value = 0.09248797987466063
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 96.12952900000001
$ws.Range("H6").Value = 288.388587
$ws.Range("I6").Value = 0.1868692395998147
$ws.Range("J6").Value = 0.1868692395998147
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 47.57896333333333
$ws.Range("N6").Value = 142.73689
$ws.Range("O6").Value = 0.450188452948237
$ws.Range("P6").Value = 0.4501884529482371
$ws.Range("Q6").Value = 4573.743335541603
$ws.Range("R6").Value = 41163.69001987443
$ws.Range("S6").Value = 0.084126373879054
$ws.Range("T6").Value = 0.08412637387905403
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 96.12952900000001
$ws.Range("H7").Value = 288.388587
$ws.Range("I7").Value = 0.1868692395998147
$ws.Range("J7").Value = 0.1868692395998147
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.829723666666666
$ws.Range("N7").Value = 29.489171
$ws.Range("O7").Value = 0.09300808131111737
$ws.Range("P7").Value = 0.09300808131111739
$ws.Range("Q7").Value = 944.9267062768197
$ws.Range("R7").Value = 8504.340356491377
$ws.Range("S7").Value = 0.01738034943124624
$ws.Range("T7").Value = 0.01738034943124624
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 96.12952900000001
$ws.Range("H8").Value = 288.388587
$ws.Range("I8").Value = 0.1868692395998147
$ws.Range("J8").Value = 0.1868692395998147
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.933664
$ws.Range("N8").Value = 32.800992
$ws.Range("O8").Value = 0.1034534789405002
$ws.Range("P8").Value = 0.1034534789405003
$ws.Range("Q8").Value = 1051.047970564256
$ws.Range("R8").Value = 9459.431735078306
$ws.Range("S8").Value = 0.01933227294356673
$ws.Range("T8").Value = 0.01933227294356673
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 96.12952900000001
$ws.Range("H9").Value = 288.388587
$ws.Range("I9").Value = 0.1868692395998147
$ws.Range("J9").Value = 0.1868692395998147
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 37.34441866666666
$ws.Range("N9").Value = 112.033256
$ws.Range("O9").Value = 0.3533499868001453
$ws.Range("P9").Value = 0.3533499868001453
$ws.Range("Q9").Value = 3589.901377205475
$ws.Range("R9").Value = 32309.11239484928
$ws.Range("S9").Value = 0.0660302433459477
$ws.Range("T9").Value = 0.06603024334594773
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 155.8267033333333
$ws.Range("H10").Value = 467.48011
$ws.Range("I10").Value = 0.3029164697274851
$ws.Range("J10").Value = 0.3029164697274851
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 47.57896333333333
$ws.Range("N10").Value = 142.73689
$ws.Range("O10").Value = 0.450188452948237
$ws.Range("P10").Value = 0.4501884529482371
$ws.Range("Q10").Value = 7414.073004250876
$ws.Range("R10").Value = 66726.65703825789
$ws.Range("S10").Value = 0.136369496879158
$ws.Range("T10").Value = 0.136369496879158
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 155.8267033333333
$ws.Range("H11").Value = 467.48011
$ws.Range("I11").Value = 0.3029164697274851
$ws.Range("J11").Value = 0.3029164697274851
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.829723666666666
$ws.Range("N11").Value = 29.489171
$ws.Range("O11").Value = 0.09300808131111737
$ws.Range("P11").Value = 0.09300808131111739
$ws.Range("Q11").Value = 1531.733433654312
$ws.Range("R11").Value = 13785.60090288881
$ws.Range("S11").Value = 0.02817367964689056
$ws.Range("T11").Value = 0.02817367964689057
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 155.8267033333333
$ws.Range("H12").Value = 467.48011
$ws.Range("I12").Value = 0.3029164697274851
$ws.Range("J12").Value = 0.3029164697274851
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 10.933664
$ws.Range("N12").Value = 32.800992
$ws.Range("O12").Value = 0.1034534789405002
$ws.Range("P12").Value = 0.1034534789405003
$ws.Range("Q12").Value = 1703.756816474347
$ws.Range("R12").Value = 15333.81134826912
$ws.Range("S12").Value = 0.03133776262168306
$ws.Range("T12").Value = 0.03133776262168307
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 155.8267033333333
$ws.Range("H13").Value = 467.48011
$ws.Range("I13").Value = 0.3029164697274851
$ws.Range("J13").Value = 0.3029164697274851
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 37.34441866666666
$ws.Range("N13").Value = 112.033256
$ws.Range("O13").Value = 0.3533499868001453
$ws.Range("P13").Value = 0.3533499868001453
$ws.Range("Q13").Value = 5819.257648726461
$ws.Range("R13").Value = 52373.31883853815
$ws.Range("S13").Value = 0.1070355305797535
$ws.Range("T13").Value = 0.1070355305797535
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 127.817347
$ws.Range("H14").Value = 383.452041
$ws.Range("I14").Value = 0.2484681938008419
$ws.Range("J14").Value = 0.248468193800842
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 47.57896333333333
$ws.Range("N14").Value = 142.73689
$ws.Range("O14").Value = 0.450188452948237
$ws.Range("P14").Value = 0.4501884529482371
$ws.Range("Q14").Value = 6081.416866276943
$ws.Range("R14").Value = 54732.75179649248
$ws.Range("S14").Value = 0.1118575117740438
$ws.Range("T14").Value = 0.1118575117740438
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 127.817347
$ws.Range("H15").Value = 383.452041
$ws.Range("I15").Value = 0.2484681938008419
$ws.Range("J15").Value = 0.248468193800842
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 9.829723666666666
$ws.Range("N15").Value = 29.489171
$ws.Range("O15").Value = 0.09300808131111737
$ws.Range("P15").Value = 0.09300808131111739
$ws.Range("Q15").Value = 1256.409200816445
$ws.Range("R15").Value = 11307.68280734801
$ws.Range("S15").Value = 0.02310954997225518
$ws.Range("T15").Value = 0.02310954997225518
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 127.817347
$ws.Range("H16").Value = 383.452041
$ws.Range("I16").Value = 0.2484681938008419
$ws.Range("J16").Value = 0.248468193800842
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 10.933664
$ws.Range("N16").Value = 32.800992
$ws.Range("O16").Value = 0.1034534789405002
$ws.Range("P16").Value = 0.1034534789405003
$ws.Range("Q16").Value = 1397.511925469408
$ws.Range("R16").Value = 12577.60732922467
$ws.Range("S16").Value = 0.02570489905475953
$ws.Range("T16").Value = 0.02570489905475954
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 127.817347
$ws.Range("H17").Value = 383.452041
$ws.Range("I17").Value = 0.2484681938008419
$ws.Range("J17").Value = 0.248468193800842
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 37.34441866666666
$ws.Range("N17").Value = 112.033256
$ws.Range("O17").Value = 0.3533499868001453
$ws.Range("P17").Value = 0.3533499868001453
$ws.Range("Q17").Value = 4773.26451923061
$ws.Range("R17").Value = 42959.38067307549
$ws.Range("S17").Value = 0.08779623299978344
$ws.Range("T17").Value = 0.08779623299978345
